$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("testA11", "testA11,11"),
    @("AAAA1", "AAAAAAAA1,"),
    @("Nathan", "Andgame1!"),
    @("lizTest", "LizTest1!"),
    @("LizTest2", "LizTest2!"),
    @("LizTester", "LizTest3!")
)

$r = 2
foreach ($pair in $data) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $r++
}
